# Generate Report for Archive
#
# 1) The localization status text changed from "Ready for handoff" to
#    "In Translation" everywhere it is used:
#      - Overview sheet: E2, F2, E3, F3 (the zh-cn / de-de status columns)
#      - zh-cn sheet:     C2, C3 (the "Status" table column)
#      - de-de sheet:     C2, C3 (the "Status" table column)
#
# 2) Because the new text is shorter, the "Status" columns were
#    auto-sized narrower:
#      - Overview sheet: columns E and F (zh-cn / de-de)
#      - zh-cn sheet:     column C (Status)
#      - de-de sheet:     column C (Status)

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet ------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
